$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New merged header "M_PL" over R1:Y1, matching the style of the existing
# "M_%cit" (B1:I1) and "M_ETR" (J1:Q1) header blocks: bold font, thin box
# border around every individual cell, centered horizontally, top-aligned
# vertically.
$cols = @("R", "S", "T", "U", "V", "W", "X", "Y")

$ws.Range("R1").Value = "M_PL"
$ws.Range("R1:Y1").Merge()

foreach ($col in $cols) {
    $cell = $ws.Range($col + "1")
    $cell.Font.Bold = $true
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Sub-headers in row 2 for the new R:Y block repeat the same sequence used
# for the other two blocks (GFA - Sales, GFA - Sales + Emp, IMF - Sales,
# IMF - Sales + Emp, OECD (20%) - Sales, OECD (20%) - Sales + Emp,
# OECD - Sales, OECD - Sales + Emp).
$subHeaders = @("GFA - Sales", "GFA - Sales + Emp", "IMF - Sales", "IMF - Sales + Emp", "OECD (20%) - Sales", "OECD (20%) - Sales + Emp", "OECD - Sales", "OECD - Sales + Emp")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "2")
    $cell.Font.Bold = $true
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Value = $subHeaders[$i]
}

# New profit data for rows 4-13, columns R:Y.
$data = @{
    4 = @(62398231596, 62571909341, 60206736130, 60380413875, 68662799453, 68662799453, 68662799453, 68662799453)
    5 = @(976531986457, 988562844368, 976531986457, 988562844368, 1017055371530, 1017055371530, 1017055371530, 1017055371530)
    6 = @(25762595315, 37793453226, 25762595315, 37793453226, 40333624448, 40333624448, 40333624448, 40333624448)
    7 = @(32205923465, 45395503555, 32205923465, 46231904709, 47942268783, 47942268783, 47942268783, 47942268783)
    8 = @(933471841988, 933471841988, 933471841988, 933471841988, 959424197928, 959424197928, 959424197928, 959424197928)
    9 = @(7937007757, 25765909108, 6548180327, 26602310262, 36502729002, 36502729002, 36502729002, 36502729002)
    10 = @(65619795685, 65619795685, 65619795685, 65619795685, 91572151625, 91572151625, 91572151625, 91572151625)
    11 = @(960129788340, 960303466085, 957938292874, 958111970619, 999081290065, 999081290065, 999081290065, 999081290065)
    12 = @(92277742037, 92451419782, 90086246571, 90259924316, 131229243762, 131229243762, 131229243762, 131229243762)
    13 = @(867852046303, 867852046303, 867852046303, 867852046303, 867852046303, 867852046303, 867852046303, 867852046303)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
